# Update "想去人数" (want-to-go count, column F) values on the 展览 and
# 全部类型 sheets to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row = new F value }
$updates = @{
    "展览" = @{
        3  = 752
        4  = 1497
        5  = 229
        8  = 6237
        11 = 115
        12 = 5195
        13 = 28
        15 = 1183
        16 = 1183
        17 = 56
        23 = 3707
    }
    "全部类型" = @{
        4  = 752
        5  = 1497
        6  = 229
        9  = 6237
        12 = 115
        13 = 5195
        14 = 28
        16 = 1183
        17 = 1183
        18 = 56
        24 = 3707
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $rowMap[$row]
    }
}
